# CET em funcao do numero de atomos
# Update existing simulation rows (2-26) with recomputed physical
# quantities and append newly-run CRN samples (crn_00025 .. crn_00029)
# as rows 27-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row = 2; B = 444.03778; C = 255.571023; D = -1.385659; E = 0.12863; F = 50; G = 0; H = 32; I = 32 },
    @{ Row = 3; B = 433.712141; C = 291.934974; D = -0.639174; E = 2.785506; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 4; B = 444.299746; C = 236.53724; D = -0.281313; E = 1.907186; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 5; B = 438.854715; C = 232.681332; D = 2.513104; E = 19.02429; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 6; B = 437.487478; C = 288.135225; D = -1.259726; E = 2.489577; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 7; B = 436.123083; C = 237.631998; D = -0.701175; E = 4.229863; F = 50; G = 1; H = 32; I = 31 },
    @{ Row = 8; B = 436.589075; C = 266.175303; D = 0.346213; E = 4.721248; F = 46.875; G = 0; H = 30; I = 34 },
    @{ Row = 9; B = 464.440069; C = 244.827381; D = -0.305409; E = 2.614112; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 10; B = 443.147867; C = 241.860944; D = -0.950894; E = 0.948705; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 11; B = 431.344354; C = 246.555231; D = 0.137237; E = 3.008565; F = 46.875; G = 0; H = 30; I = 34 },
    @{ Row = 12; B = 435.123183; C = 231.730715; D = 0.703753; E = 3.837027; F = 50; G = 0; H = 32; I = 32 },
    @{ Row = 13; B = 468.192422; C = 203.046293; D = 0.598958; E = 2.809231; F = 53.125; G = 0; H = 34; I = 30 },
    @{ Row = 14; B = 408.617598; C = 245.175361; D = -0.75159; E = 2.282665; F = 53.125; G = 1; H = 34; I = 29 },
    @{ Row = 15; B = 433.729331; C = 246.581486; D = 0.04456; E = 4.132562; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 16; B = 423.564751; C = 290.307411; D = 0.205954; E = 5.051874; F = 50; G = 1; H = 32; I = 31 },
    @{ Row = 17; B = 429.975678; C = 257.803029; D = -0.508456; E = 1.116876; F = 50; G = 1; H = 32; I = 31 },
    @{ Row = 18; B = 436.986294; C = 271.552346; D = 0.283232; E = 4.043613; F = 50; G = 1; H = 32; I = 31 },
    @{ Row = 19; B = 420.981257; C = 274.686059; D = -1.327174; E = -0.942836; F = 50; G = 1; H = 32; I = 31 },
    @{ Row = 20; B = 445.751152; C = 249.564427; D = 1.837633; E = 11.44626; F = 50; G = 0; H = 32; I = 32 },
    @{ Row = 21; B = 431.464251; C = 249.325769; D = -0.839351; E = 2.392735; F = 53.125; G = 1; H = 34; I = 29 },
    @{ Row = 22; B = 448.065918; C = 268.772782; D = 1.058055; E = 7.35721; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 23; B = 431.858355; C = 251.761421; D = -0.400146; E = 3.743755; F = 53.125; G = 1; H = 34; I = 29 },
    @{ Row = 24; B = 435.054281; C = 285.61565; D = 0.575845; E = 4.082143; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 25; B = 432.287002; C = 234.665684; D = 1.43741; E = 4.910842; F = 46.875; G = 1; H = 30; I = 33 },
    @{ Row = 26; B = 435.657422; C = 270.027727; D = 0.305117; E = 14.644; F = 53.125; G = 0; H = 34; I = 30 },
    @{ Row = 27; B = 436.485493; C = 273.333769; D = 0.673169; E = 5.818353; F = 50; G = 0; H = 32; I = 32 },
    @{ Row = 28; B = 436.520019; C = 251.186596; D = 0.778664; E = 3.147818; F = 53.125; G = 0; H = 34; I = 30 },
    @{ Row = 29; B = 437.885241; C = 247.727518; D = 0.528135; E = 6.927249; F = 50; G = 1; H = 32; I = 31 },
    @{ Row = 30; B = 443.374594; C = 284.04816; D = 0.124317; E = 16.030135; F = 50; G = 0; H = 32; I = 32 },
    @{ Row = 31; B = 453.126405; C = 233.209148; D = -0.948704; E = 1.171624; F = 46.875; G = 0; H = 30; I = 34 }
)

$newLabels = @{
    27 = "crn_00025"
    28 = "crn_00026"
    29 = "crn_00027"
    30 = "crn_00028"
    31 = "crn_00029"
}

foreach ($entry in $rowData) {
    $r = $entry.Row

    if ($newLabels.ContainsKey($r)) {
        # Brand-new row: write the CRN label first, then give the whole
        # row the same border/format style used by the existing data rows.
        $ws.Cells.Item($r, 1).Value = $newLabels[$r]
        $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 9)).Borders.LineStyle = 1
    }

    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
}
